# Move bus speed factors from hardcoded to additional_inputs.xlsx (Buskosten sheet)
# Adds rows 53-59 below the existing "Gemiddelde snelheid per bustype" block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Buskosten")

# Row 53 left blank (spacer row between the two variable blocks)

# Row 54: section header
$ws.Cells.Item(54, 1).Value = "Snelheidsfactor bus t.o.v. auto (Google Maps geeft autosnelheid):"

# Rows 55-59: new configurable speed factor variables
$ws.Cells.Item(55, 1).Value = "snelheidsfactor_touringcar"
$ws.Cells.Item(55, 2).Value = 0.95

$ws.Cells.Item(56, 1).Value = "snelheidsfactor_dubbeldekker"
$ws.Cells.Item(56, 2).Value = 0.9

$ws.Cells.Item(57, 1).Value = "snelheidsfactor_lagevloer"
$ws.Cells.Item(57, 2).Value = 0.85

$ws.Cells.Item(58, 1).Value = "snelheidsfactor_midibus"
$ws.Cells.Item(58, 2).Value = 0.92

$ws.Cells.Item(59, 1).Value = "snelheidsfactor_taxibus"
$ws.Cells.Item(59, 2).Value = 0.95
